$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$s.Delete()
